# Re-generate statistics: fix minutes/seconds formatting in the "haul" (total time)
# column so single-digit minutes/seconds are zero-padded, e.g.
#   "4 ч. 7 мин. 54 сек."  ->  "4 ч. 07 мин. 54 сек."
#   "2 ч. 26 мин. 3 сек."  ->  "2 ч. 26 мин. 03 сек."
# The hours part is left untouched (can stay single- or multi-digit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

# Matches "<hours> ч. <minutes> мин. <seconds> сек." strings.
$pattern = [regex]'^(\d+) ч\. (\d+) мин\. (\d+) сек\.$'

$changed = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string]) {
            $m = $pattern.Match($val)
            if ($m.Success) {
                $hours = $m.Groups[1].Value
                $minutes = $m.Groups[2].Value.PadLeft(2, '0')
                $seconds = $m.Groups[3].Value.PadLeft(2, '0')
                $newVal = "$hours ч. $minutes мин. $seconds сек."
                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                    $changed = $changed + 1
                }
            }
        }
    }
}

Write-Output "Re-formatted $changed haul-time cell(s)."
